$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28 (Cause=Cardiorespiratory, Sex=1, Age=55)
$ws.Range("D28").Value = 120523.83
$ws.Range("E28").Value = 120528.28
$ws.Range("F28").Formula = "=IF(D28<E28,""non-pw"",""pw"")"

# Row 29 (Cause=Cardiorespiratory, Sex=1, Age=65)
$ws.Range("D29").Value = 129515.48
$ws.Range("E29").Value = 129516.42
$ws.Range("F29").Formula = "=IF(D29<E29,""non-pw"",""pw"")"

# Row 38 (Cause=Cardiorespiratory, Sex=2, Age=55)
$ws.Range("D38").Value = 113458.35
$ws.Range("E38").Value = 113463.93
$ws.Range("F38").Formula = "=IF(D38<E38,""non-pw"",""pw"")"

# Row 39 (Cause=Cardiorespiratory, Sex=2, Age=65)
$ws.Range("D39").Value = 123525.49
$ws.Range("E39").Value = 123530.41
$ws.Range("F39").Formula = "=IF(D39<E39,""non-pw"",""pw"")"

# Update sheet view (top left cell + selection) to match the target state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F29").Select()
